$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.674.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.796.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.70%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  +4.19%  '
$ws.Range('E9').Value = '  +2.40%  '
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0951'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.057.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.798.33'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.590.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0812'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.64%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '248.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('E26').Value = '  +1.42%  '
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('E28').Value = '  +2.82%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +10.00%  '
$ws.Range('E31').Value = '  +3.73%  '
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +2.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.422.76'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  +6.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.675'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0193'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.06'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.72%  '
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.937'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.57'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('E45').Value = '  +3.08%  '
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.957.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.28%  '
